# Scheduled-runner refresh of market price / profit figures across the
# per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose underlying market data changed since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2602.1333
$ws.Range("J17").Value = 2670.3333
$ws.Range("L17").Value = 8010.999899999999
$ws.Range("N17").Value = -8346.999899999999
$ws.Range("H53").Value = 2668.6924
$ws.Range("I53").Value = 2442.7144
$ws.Range("J53").Value = 2932.3333
$ws.Range("K53").Value = 2442.7144
$ws.Range("L53").Value = 2932.3333
$ws.Range("M53").Value = -1805.7144
$ws.Range("N53").Value = -4206.3333
$ws.Range("H137").Value = 970173.1
$ws.Range("I137").Value = 2914.889
$ws.Range("K137").Value = 8744.667000000001
$ws.Range("M137").Value = -6194.667000000001
$ws.Range("H138").Value = 1833.8909
$ws.Range("I138").Value = 1191.9032
$ws.Range("J138").Value = 2663.125
$ws.Range("K138").Value = 3575.7096
$ws.Range("L138").Value = 7989.375
$ws.Range("M138").Value = 1564.2904
$ws.Range("N138").Value = -18269.375
$ws.Range("H141").Value = 884.125
$ws.Range("I141").Value = 883.0909
$ws.Range("J141").Value = 895.5
$ws.Range("K141").Value = 2649.2727
$ws.Range("L141").Value = 2686.5
$ws.Range("M141").Value = 2530.7273
$ws.Range("N141").Value = -13046.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1083.2222
$ws.Range("I2").Value = 960.46155
$ws.Range("J2").Value = 1402.4
$ws.Range("K2").Value = 960.46155
$ws.Range("L2").Value = 1402.4
$ws.Range("M2").Value = -847.46155
$ws.Range("N2").Value = -1628.4
$ws.Range("H32").Value = 4613.8
$ws.Range("I32").Value = 3228.2673
$ws.Range("J32").Value = 13124.929
$ws.Range("K32").Value = 3228.2673
$ws.Range("L32").Value = 13124.929
$ws.Range("M32").Value = -2941.2673
$ws.Range("N32").Value = -13698.929
$ws.Range("H45").Value = 5441354.5
$ws.Range("I45").Value = 7759.778
$ws.Range("J45").Value = 25002296
$ws.Range("K45").Value = 7759.778
$ws.Range("L45").Value = 25002296
$ws.Range("M45").Value = -7382.778
$ws.Range("N45").Value = -25003050
$ws.Range("H110").Value = 831.73334
$ws.Range("I110").Value = 813.9091
$ws.Range("J110").Value = 880.75
$ws.Range("K110").Value = 813.9091
$ws.Range("L110").Value = 880.75
$ws.Range("M110").Value = 1231.0909
$ws.Range("N110").Value = -4970.75
$ws.Range("H116").Value = 1083.2222
$ws.Range("I116").Value = 960.46155
$ws.Range("J116").Value = 1402.4
$ws.Range("K116").Value = 960.46155
$ws.Range("L116").Value = 1402.4
$ws.Range("M116").Value = 1333.53845
$ws.Range("N116").Value = -5990.4
$ws.Range("H122").Value = 3191.0286
$ws.Range("I122").Value = 3926.3
$ws.Range("J122").Value = 2210.6667
$ws.Range("K122").Value = 11778.9
$ws.Range("L122").Value = 6632.000100000001
$ws.Range("M122").Value = -9328.900000000001
$ws.Range("N122").Value = -11532.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1083.2222
$ws.Range("I3").Value = 960.46155
$ws.Range("J3").Value = 1402.4
$ws.Range("K3").Value = 960.46155
$ws.Range("L3").Value = 1402.4
$ws.Range("M3").Value = -846.46155
$ws.Range("N3").Value = -1630.4
$ws.Range("H94").Value = 1462.5946
$ws.Range("I94").Value = 1565.4193
$ws.Range("J94").Value = 931.3333
$ws.Range("K94").Value = 1565.4193
$ws.Range("L94").Value = 931.3333
$ws.Range("M94").Value = -1114.4193
$ws.Range("N94").Value = -1833.3333
$ws.Range("H134").Value = 1975.4706
$ws.Range("I134").Value = 1129.4615
$ws.Range("J134").Value = 4725
$ws.Range("K134").Value = 3388.3845
$ws.Range("L134").Value = 14175
$ws.Range("M134").Value = -853.3844999999997
$ws.Range("N134").Value = -19245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3473.087
$ws.Range("I31").Value = 2185
$ws.Range("J31").Value = 4653.8335
$ws.Range("K31").Value = 2185
$ws.Range("L31").Value = 4653.8335
$ws.Range("M31").Value = -1890
$ws.Range("N31").Value = -5243.8335
$ws.Range("H34").Value = 3473.087
$ws.Range("I34").Value = 2185
$ws.Range("J34").Value = 4653.8335
$ws.Range("K34").Value = 2185
$ws.Range("L34").Value = 4653.8335
$ws.Range("M34").Value = -1983
$ws.Range("N34").Value = -5057.8335
$ws.Range("H99").Value = 1491553.8
$ws.Range("I99").Value = 3260
$ws.Range("J99").Value = 2844548
$ws.Range("K99").Value = 3260
$ws.Range("L99").Value = 2844548
$ws.Range("M99").Value = -1762
$ws.Range("N99").Value = -2847544
$ws.Range("H126").Value = 1491553.8
$ws.Range("I126").Value = 3260
$ws.Range("J126").Value = 2844548
$ws.Range("K126").Value = 9780
$ws.Range("L126").Value = 8533644
$ws.Range("M126").Value = -7310
$ws.Range("N126").Value = -8538584
$ws.Range("H132").Value = 1803.1904
$ws.Range("I132").Value = 1483.3529
$ws.Range("J132").Value = 3162.5
$ws.Range("K132").Value = 4450.0587
$ws.Range("L132").Value = 9487.5
$ws.Range("M132").Value = -1920.0587
$ws.Range("N132").Value = -14547.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 899.25
$ws.Range("I64").Value = 899.25
$ws.Range("K64").Value = 2697.75
$ws.Range("M64").Value = -2427.75
$ws.Range("H67").Value = 899.25
$ws.Range("I67").Value = 899.25
$ws.Range("K67").Value = 2697.75
$ws.Range("M67").Value = -1761.75
$ws.Range("H114").Value = 7338.4
$ws.Range("I114").Value = 999.5
$ws.Range("J114").Value = 8923.125
$ws.Range("K114").Value = 2998.5
$ws.Range("L114").Value = 26769.375
$ws.Range("M114").Value = 255.5
$ws.Range("N114").Value = -33277.375
$ws.Range("H137").Value = 4577.3447
$ws.Range("I137").Value = 3574.875
$ws.Range("J137").Value = 4959.2383
$ws.Range("K137").Value = 10724.625
$ws.Range("L137").Value = 14877.7149
$ws.Range("M137").Value = -5624.625
$ws.Range("N137").Value = -25077.7149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2186.5
$ws.Range("I80").Value = 999
$ws.Range("J80").Value = 3374
$ws.Range("K80").Value = 999
$ws.Range("L80").Value = 3374
$ws.Range("M80").Value = -1
$ws.Range("N80").Value = -5370
$ws.Range("H83").Value = 2186.5
$ws.Range("I83").Value = 999
$ws.Range("J83").Value = 3374
$ws.Range("K83").Value = 4995
$ws.Range("L83").Value = 16870
$ws.Range("M83").Value = -3
$ws.Range("N83").Value = -26854

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 9990
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H136").Value = 4003.5454
$ws.Range("I136").Value = 5120.364
$ws.Range("J136").Value = 3445.1365
$ws.Range("K136").Value = 15361.092
$ws.Range("L136").Value = 10335.4095
$ws.Range("M136").Value = -12811.092
$ws.Range("N136").Value = -15435.4095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 763719.9399999999
$ws.Range("I132").Value = 753.7659
$ws.Range("J132").Value = 4349661
$ws.Range("K132").Value = 2261.2977
$ws.Range("L132").Value = 13048983
$ws.Range("M132").Value = 268.7022999999999
$ws.Range("N132").Value = -13054043
